# Community data update per commit: "minor bug fix on community entity"
# Replaces rows 2-17 with new barangay data and appends rows 18-30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $val) {
    # Leading apostrophe forces text storage (matches source data which
    # stores every value, numeric-looking or not, as text).
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

Set-TextCell $ws 2 1 "Balite"
Set-TextCell $ws 2 2 "14.8956"
Set-TextCell $ws 2 3 "120.7855"
Set-TextCell $ws 2 4 "5016"
Set-TextCell $ws 2 5 "0"
Set-TextCell $ws 2 6 "2144"
Set-TextCell $ws 2 7 "1000"
Set-TextCell $ws 2 8 ""

Set-TextCell $ws 3 1 "Balungao"
Set-TextCell $ws 3 2 "14.9143"
Set-TextCell $ws 3 3 "120.7622"
Set-TextCell $ws 3 4 "5720"
Set-TextCell $ws 3 5 "0"
Set-TextCell $ws 3 6 "3366"
Set-TextCell $ws 3 7 "1000"
Set-TextCell $ws 3 8 ""

Set-TextCell $ws 4 1 "Buguion"
Set-TextCell $ws 4 2 "14.894"
Set-TextCell $ws 4 3 "120.7985"
Set-TextCell $ws 4 4 "3841"
Set-TextCell $ws 4 5 "0"
Set-TextCell $ws 4 6 "2196"
Set-TextCell $ws 4 7 "1000"
Set-TextCell $ws 4 8 ""

Set-TextCell $ws 5 1 "Bulusan"
Set-TextCell $ws 5 2 "14.9076"
Set-TextCell $ws 5 3 "120.7455"
Set-TextCell $ws 5 4 "2603"
Set-TextCell $ws 5 5 "0"
Set-TextCell $ws 5 6 "1721"
Set-TextCell $ws 5 7 "1000"
Set-TextCell $ws 5 8 ""

Set-TextCell $ws 6 1 "Calizon"
Set-TextCell $ws 6 2 "14.9125"
Set-TextCell $ws 6 3 "120.753"
Set-TextCell $ws 6 4 "2221"
Set-TextCell $ws 6 5 "0"
Set-TextCell $ws 6 6 "1387"
Set-TextCell $ws 6 7 "1000"
Set-TextCell $ws 6 8 ""

Set-TextCell $ws 7 1 "Calumpang"
Set-TextCell $ws 7 2 "14.8845"
Set-TextCell $ws 7 3 "120.7838"
Set-TextCell $ws 7 4 "3517"
Set-TextCell $ws 7 5 "0"
Set-TextCell $ws 7 6 "2784"
Set-TextCell $ws 7 7 "1000"
Set-TextCell $ws 7 8 ""

Set-TextCell $ws 8 1 "Caniogan"
Set-TextCell $ws 8 2 "14.9054"
Set-TextCell $ws 8 3 "120.7733"
Set-TextCell $ws 8 4 "4510"
Set-TextCell $ws 8 5 "0"
Set-TextCell $ws 8 6 "2869"
Set-TextCell $ws 8 7 "1000"
Set-TextCell $ws 8 8 ""

Set-TextCell $ws 9 1 "Corazon"
Set-TextCell $ws 9 2 "14.9128"
Set-TextCell $ws 9 3 "120.7686"
Set-TextCell $ws 9 4 "2175"
Set-TextCell $ws 9 5 "0"
Set-TextCell $ws 9 6 "1647"
Set-TextCell $ws 9 7 "1000"
Set-TextCell $ws 9 8 ""

Set-TextCell $ws 10 1 "Frances"
Set-TextCell $ws 10 2 "14.9153"
Set-TextCell $ws 10 3 "120.7532"
Set-TextCell $ws 10 4 "6129"
Set-TextCell $ws 10 5 "6"
Set-TextCell $ws 10 6 "3819"
Set-TextCell $ws 10 7 "1000"
Set-TextCell $ws 10 8 ""

Set-TextCell $ws 11 1 "Gatbuca"
Set-TextCell $ws 11 2 "14.9218"
Set-TextCell $ws 11 3 "120.7685"
Set-TextCell $ws 11 4 "6384"
Set-TextCell $ws 11 5 "115"
Set-TextCell $ws 11 6 "4250"
Set-TextCell $ws 11 7 "1000"
Set-TextCell $ws 11 8 ""

Set-TextCell $ws 12 1 "Gugo"
Set-TextCell $ws 12 2 "14.9014"
Set-TextCell $ws 12 3 "120.7548"
Set-TextCell $ws 12 4 "1960"
Set-TextCell $ws 12 5 "57"
Set-TextCell $ws 12 6 "1179"
Set-TextCell $ws 12 7 "1000"
Set-TextCell $ws 12 8 ""

Set-TextCell $ws 13 1 "Iba Este"
Set-TextCell $ws 13 2 "14.8899"
Set-TextCell $ws 13 3 "120.7673"
Set-TextCell $ws 13 4 "4161"
Set-TextCell $ws 13 5 "0"
Set-TextCell $ws 13 6 "1828"
Set-TextCell $ws 13 7 "1000"
Set-TextCell $ws 13 8 ""

Set-TextCell $ws 14 1 "Iba O'Este"
Set-TextCell $ws 14 2 "14.8919"
Set-TextCell $ws 14 3 "120.7635"
Set-TextCell $ws 14 4 "14085"
Set-TextCell $ws 14 5 "601"
Set-TextCell $ws 14 6 "8095"
Set-TextCell $ws 14 7 "1000"
Set-TextCell $ws 14 8 ""

Set-TextCell $ws 15 1 "Longos"
Set-TextCell $ws 15 2 "14.8748"
Set-TextCell $ws 15 3 "120.7866"
Set-TextCell $ws 15 4 "4265"
Set-TextCell $ws 15 5 "0"
Set-TextCell $ws 15 6 "2293"
Set-TextCell $ws 15 7 "1000"
Set-TextCell $ws 15 8 ""

Set-TextCell $ws 16 1 "Meysulao"
Set-TextCell $ws 16 2 "14.9078"
Set-TextCell $ws 16 3 "120.7397"
Set-TextCell $ws 16 4 "4280"
Set-TextCell $ws 16 5 "56"
Set-TextCell $ws 16 6 "2687"
Set-TextCell $ws 16 7 "1000"
Set-TextCell $ws 16 8 ""

Set-TextCell $ws 17 1 "Meyto"
Set-TextCell $ws 17 2 "14.8831"
Set-TextCell $ws 17 3 "120.7295"
Set-TextCell $ws 17 4 "2925"
Set-TextCell $ws 17 5 "6"
Set-TextCell $ws 17 6 "1975"
Set-TextCell $ws 17 7 "1000"
Set-TextCell $ws 17 8 ""

Set-TextCell $ws 18 1 "Palimbang"
Set-TextCell $ws 18 2 "14.8994"
Set-TextCell $ws 18 3 "120.7756"
Set-TextCell $ws 18 4 "1684"
Set-TextCell $ws 18 5 "0"
Set-TextCell $ws 18 6 "1424"
Set-TextCell $ws 18 7 "1000"
Set-TextCell $ws 18 8 ""

Set-TextCell $ws 19 1 "Panducot"
Set-TextCell $ws 19 2 "14.8761"
Set-TextCell $ws 19 3 "120.738"
Set-TextCell $ws 19 4 "1752"
Set-TextCell $ws 19 5 "0"
Set-TextCell $ws 19 6 "1713"
Set-TextCell $ws 19 7 "1000"
Set-TextCell $ws 19 8 ""

Set-TextCell $ws 20 1 "Pio Cruzcosa"
Set-TextCell $ws 20 2 "14.8881"
Set-TextCell $ws 20 3 "120.7855"
Set-TextCell $ws 20 4 "4663"
Set-TextCell $ws 20 5 "92"
Set-TextCell $ws 20 6 "2899"
Set-TextCell $ws 20 7 "1000"
Set-TextCell $ws 20 8 ""

Set-TextCell $ws 21 1 "Poblacion"
Set-TextCell $ws 21 2 "14.9157"
Set-TextCell $ws 21 3 "120.7672"
Set-TextCell $ws 21 4 "1785"
Set-TextCell $ws 21 5 "0"
Set-TextCell $ws 21 6 "1294"
Set-TextCell $ws 21 7 "1000"
Set-TextCell $ws 21 8 ""

Set-TextCell $ws 22 1 "Pungo"
Set-TextCell $ws 22 2 "14.9023"
Set-TextCell $ws 22 3 "120.7914"
Set-TextCell $ws 22 4 "9528"
Set-TextCell $ws 22 5 "0"
Set-TextCell $ws 22 6 "5486"
Set-TextCell $ws 22 7 "1000"
Set-TextCell $ws 22 8 ""

Set-TextCell $ws 23 1 "San Jose"
Set-TextCell $ws 23 2 "14.8838"
Set-TextCell $ws 23 3 "120.7395"
Set-TextCell $ws 23 4 "5661"
Set-TextCell $ws 23 5 "0"
Set-TextCell $ws 23 6 "3629"
Set-TextCell $ws 23 7 "1000"
Set-TextCell $ws 23 8 ""

Set-TextCell $ws 24 1 "San Marcos"
Set-TextCell $ws 24 2 "14.8976"
Set-TextCell $ws 24 3 "120.7797"
Set-TextCell $ws 24 4 "2671"
Set-TextCell $ws 24 5 "0"
Set-TextCell $ws 24 6 "1471"
Set-TextCell $ws 24 7 "1000"
Set-TextCell $ws 24 8 ""

Set-TextCell $ws 25 1 "San Miguel"
Set-TextCell $ws 25 2 "14.917"
Set-TextCell $ws 25 3 "120.7427"
Set-TextCell $ws 25 4 "6005"
Set-TextCell $ws 25 5 "17"
Set-TextCell $ws 25 6 "3287"
Set-TextCell $ws 25 7 "1000"
Set-TextCell $ws 25 8 ""

Set-TextCell $ws 26 1 "Santa Lucia"
Set-TextCell $ws 26 2 "14.8982"
Set-TextCell $ws 26 3 "120.736"
Set-TextCell $ws 26 4 "2460"
Set-TextCell $ws 26 5 "0"
Set-TextCell $ws 26 6 "1483"
Set-TextCell $ws 26 7 "1000"
Set-TextCell $ws 26 8 ""

Set-TextCell $ws 27 1 "Santo Niño"
Set-TextCell $ws 27 2 "14.9047"
Set-TextCell $ws 27 3 "120.7792"
Set-TextCell $ws 27 4 "2544"
Set-TextCell $ws 27 5 "0"
Set-TextCell $ws 27 6 "1392"
Set-TextCell $ws 27 7 "1000"
Set-TextCell $ws 27 8 ""

Set-TextCell $ws 28 1 "Sapang Bayan"
Set-TextCell $ws 28 2 "14.9196"
Set-TextCell $ws 28 3 "120.7739"
Set-TextCell $ws 28 4 "3140"
Set-TextCell $ws 28 5 "65"
Set-TextCell $ws 28 6 "1775"
Set-TextCell $ws 28 7 "1000"
Set-TextCell $ws 28 8 ""

Set-TextCell $ws 29 1 "Sergio Bayan"
Set-TextCell $ws 29 2 "14.894"
Set-TextCell $ws 29 3 "120.7909"
Set-TextCell $ws 29 4 "1727"
Set-TextCell $ws 29 5 "0"
Set-TextCell $ws 29 6 "1258"
Set-TextCell $ws 29 7 "1000"
Set-TextCell $ws 29 8 ""

Set-TextCell $ws 30 1 "Sucol"
Set-TextCell $ws 30 2 "14.9138"
Set-TextCell $ws 30 3 "120.7701"
Set-TextCell $ws 30 4 "1059"
Set-TextCell $ws 30 5 "0"
Set-TextCell $ws 30 6 "963"
Set-TextCell $ws 30 7 "1000"
Set-TextCell $ws 30 8 ""
